$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("D11").Value = 44505
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 19000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 19500
$ws.Range("Q11").Value = '$/bandeja 8 kilos'
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 2438
$ws.Range("T11").Value = 8

# Row 12
$ws.Range("D12").Value = 44873
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 22000
$ws.Range("O12").Value = 22500
$ws.Range("P12").Value = 22250
$ws.Range("Q12").Value = '$/bandeja 8 kilos'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 2781
$ws.Range("T12").Value = 8

# Row 13
$ws.Range("D13").Value = 44488
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 160
$ws.Range("N13").Value = 17000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 17500
$ws.Range("Q13").Value = '$/bandeja 8 kilos'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 2188
$ws.Range("T13").Value = 8

# Row 14
$ws.Range("D14").Value = 44895
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 22000
$ws.Range("O14").Value = 22500
$ws.Range("P14").Value = 22250
$ws.Range("Q14").Value = '$/bandeja 8 kilos'
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 2781
$ws.Range("T14").Value = 8

# Row 15
$ws.Range("D15").Value = 44519
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 400
$ws.Range("N15").Value = 21000
$ws.Range("O15").Value = 22500
$ws.Range("P15").Value = 21500
$ws.Range("Q15").Value = '$/bandeja 8 kilos'
$ws.Range("R15").Value = 'Provincia de Limarí'
$ws.Range("S15").Value = 2688
$ws.Range("T15").Value = 8

# Row 16
$ws.Range("D16").Value = 44519
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 18000
$ws.Range("Q16").Value = '$/bandeja 8 kilos'
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 2250
$ws.Range("T16").Value = 8

# Row 17
$ws.Range("D17").Value = 44890
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 22000
$ws.Range("O17").Value = 22500
$ws.Range("P17").Value = 22250
$ws.Range("Q17").Value = '$/bandeja 8 kilos'
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 2781
$ws.Range("T17").Value = 8

# Row 18
$ws.Range("D18").Value = 45240
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 150
$ws.Range("N18").Value = 28000
$ws.Range("O18").Value = 28000
$ws.Range("P18").Value = 28000
$ws.Range("Q18").Value = '$/bandeja 10 kilos'
$ws.Range("R18").Value = 'Provincia del Elquí'
$ws.Range("S18").Value = 2800
$ws.Range("T18").Value = 10

# Row 19
$ws.Range("D19").Value = 44533
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 19000
$ws.Range("P19").Value = 18500
$ws.Range("Q19").Value = '$/bandeja 8 kilos'
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 2312
$ws.Range("T19").Value = 8

# Row 20
$ws.Range("D20").Value = 44533
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 16000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 16000
$ws.Range("Q20").Value = '$/bandeja 8 kilos'
$ws.Range("R20").Value = 'Provincia de Limarí'
$ws.Range("S20").Value = 2000
$ws.Range("T20").Value = 8

# Row 21
$ws.Range("D21").Value = 44516
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 19000
$ws.Range("P21").Value = 18500
$ws.Range("Q21").Value = '$/bandeja 8 kilos'
$ws.Range("R21").Value = 'Provincia de Limarí'
$ws.Range("S21").Value = 2312
$ws.Range("T21").Value = 8

# Row 22
$ws.Range("D22").Value = 44876
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 22000
$ws.Range("O22").Value = 22500
$ws.Range("P22").Value = 22250
$ws.Range("Q22").Value = '$/bandeja 8 kilos'
$ws.Range("R22").Value = 'Provincia de Limarí'
$ws.Range("S22").Value = 2781
$ws.Range("T22").Value = 8

# Row 23
$ws.Range("D23").Value = 44159
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 300
$ws.Range("N23").Value = 2000
$ws.Range("O23").Value = 2100
$ws.Range("P23").Value = 2050
$ws.Range("Q23").Value = '$/kilo (en caja de 14 kilos)'
$ws.Range("R23").Value = 'Provincia de Limarí'
$ws.Range("S23").Value = 2050
$ws.Range("T23").Value = 1

# Row 24
$ws.Range("D24").Value = 44880
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 300
$ws.Range("N24").Value = 22000
$ws.Range("O24").Value = 22500
$ws.Range("P24").Value = 22250
$ws.Range("Q24").Value = '$/bandeja 8 kilos'
$ws.Range("R24").Value = 'Provincia de Limarí'
$ws.Range("S24").Value = 2781
$ws.Range("T24").Value = 8

# Row 25
$ws.Range("D25").Value = 44526
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 21000
$ws.Range("O25").Value = 21000
$ws.Range("P25").Value = 21000
$ws.Range("Q25").Value = '$/bandeja 8 kilos'
$ws.Range("R25").Value = 'Provincia de Limarí'
$ws.Range("S25").Value = 2625
$ws.Range("T25").Value = 8

# Row 26
$ws.Range("D26").Value = 45244
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 150
$ws.Range("N26").Value = 25000
$ws.Range("O26").Value = 25000
$ws.Range("P26").Value = 25000
$ws.Range("Q26").Value = '$/bandeja 10 kilos'
$ws.Range("R26").Value = 'Provincia del Elquí'
$ws.Range("S26").Value = 2500
$ws.Range("T26").Value = 10
